$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

# Add a new row at the end of the table (clones formatting from the last row)
$newRow = $table.Rows.Add()
$rowIndex = $newRow.Index

# Clear the auto-populated placeholder runs in columns 1-3 so they come back
# to plain empty paragraphs (column 3 keeps its paragraph-mark formatting).
$table.Cell($rowIndex, 1).Range.Delete()
$table.Cell($rowIndex, 2).Range.Delete()
$table.Cell($rowIndex, 3).Range.Delete()

# Column 4 gets the new name.
$table.Cell($rowIndex, 4).Range.Delete()
$table.Cell($rowIndex, 4).Range.Text = "Keramat Ali"
